$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.447.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.853.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "460.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.75%  "

$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.14%  "

$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000315"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.467.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("E15").Value = "  -5.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.863.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("E19").Value = "  +7.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.576.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +7.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "752.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "

$ws.Range("E31").Value = "  +11.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.20%  "

$ws.Range("E35").Value = "  +7.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("E39").Value = "  +4.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.356"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0675"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.73%  "

$ws.Range("E44").Value = "  +5.18%  "

$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("E48").Value = "  +7.75%  "

$ws.Range("E49").Value = "  +4.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.87%  "
